$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 & 4 mirror the existing data layout: columns A-C are text
# (date/time strings stored as text, like row 2), D/E/H are numbers,
# F/G are text looked up from the existing shared strings.
$ws.Range("A3:C4").NumberFormat = "@"

# --- Row 3: "Physik Lerngruppe" on 2025-04-20, 08:00-09:00 ---
$ws.Range("A3").Value = "2025-04-20"
$ws.Range("B3").Value = "08:00"
$ws.Range("C3").Value = "09:00"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Physik Lerngruppe"
$ws.Range("G3").Value = "Erste Verantstaltung"
$ws.Range("H3").Value = 5

# --- Row 4: "Physik Lerngruppe" on 2025-04-16, 11:00-14:00 ---
# (B4/C4 before A4 so new shared-string entries land in the same order
# Excel produced them in.)
$ws.Range("B4").Value = "11:00"
$ws.Range("C4").Value = "14:00"
$ws.Range("A4").Value = "2025-04-16"
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Physik Lerngruppe"
$ws.Range("G4").Value = "Erste Verantstaltung"
$ws.Range("H4").Value = 6

[void]$ws.Range("A4").Select()
